$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colG = 7

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
